# Auto-generated: refresh static market-price columns (H-N) across sheets
# per scheduled-runner data pull. Source cell values/removals enumerated
# directly from the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 428
$ws.Range("I2").Value = 439.6
$ws.Range("K2").Value = 439.6
$ws.Range("M2").Value = -326.6

$ws.Range("H80").Value = 729.32355
$ws.Range("J80").Value = 985.6667
$ws.Range("L80").Value = 2957.0001
$ws.Range("N80").Value = -4953.0001

$ws.Range("H83").Value = 729.32355
$ws.Range("J83").Value = 985.6667
$ws.Range("L83").Value = 8871.0003
$ws.Range("N83").Value = -18855.0003

$ws.Range("H86").Value = 2097.125
$ws.Range("I86").Value = 1500
$ws.Range("K86").Value = 1500
$ws.Range("M86").Value = -377

$ws.Range("H89").Value = 2097.125
$ws.Range("I89").Value = 1500
$ws.Range("K89").Value = 7500
$ws.Range("M89").Value = -1884

$ws.Range("H116").Value = 5493.6875
$ws.Range("I116").Value = 5175.25
$ws.Range("J116").Value = 5812.125
$ws.Range("K116").Value = 5175.25
$ws.Range("L116").Value = 5812.125
$ws.Range("M116").Value = -1733.25
$ws.Range("N116").Value = -12696.125

$ws.Range("H124").Value = 100500
$ws.Range("J124").Value = 100500
$ws.Range("L124").Value = 100500
$ws.Range("N124").Value = -110320

$ws.Range("H132").Value = 1102.2727
$ws.Range("I132").Value = 762.5
$ws.Range("K132").Value = 2287.5
$ws.Range("M132").Value = 242.5

$ws.Range("I137").Value = 2108.85
$ws.Range("J137").Value = 3754
$ws.Range("K137").Value = 6326.549999999999
$ws.Range("L137").Value = 11262
$ws.Range("M137").Value = -3776.549999999999
$ws.Range("N137").Value = -16362

$ws.Range("H141").Value = 4315
$ws.Range("I141").Value = 2522.5
$ws.Range("K141").Value = 7567.5
$ws.Range("M141").Value = -2387.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11364383
$ws.Range("I2").Value = 22727974
$ws.Range("J2").Value = 791.75
$ws.Range("K2").Value = 22727974
$ws.Range("L2").Value = 791.75
$ws.Range("M2").Value = -22727861
$ws.Range("N2").Value = -1017.75

$ws.Range("H74").Value = 2885.0833
$ws.Range("I74").Value = 1151.5714
$ws.Range("J74").Value = 5312
$ws.Range("K74").Value = 1151.5714
$ws.Range("L74").Value = 5312
$ws.Range("M74").Value = -277.5714
$ws.Range("N74").Value = -7060

$ws.Range("H77").Value = 2885.0833
$ws.Range("I77").Value = 1151.5714
$ws.Range("J77").Value = 5312
$ws.Range("K77").Value = 5757.857
$ws.Range("L77").Value = 26560
$ws.Range("M77").Value = -1389.857
$ws.Range("N77").Value = -35296

$ws.Range("H116").Value = 11364383
$ws.Range("I116").Value = 22727974
$ws.Range("J116").Value = 791.75
$ws.Range("K116").Value = 22727974
$ws.Range("L116").Value = 791.75
$ws.Range("M116").Value = -22725680
$ws.Range("N116").Value = -5379.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11364383
$ws.Range("I3").Value = 22727974
$ws.Range("J3").Value = 791.75
$ws.Range("K3").Value = 22727974
$ws.Range("L3").Value = 791.75
$ws.Range("M3").Value = -22727860
$ws.Range("N3").Value = -1019.75

$ws.Range("H11").Value = 437
$ws.Range("I11").Value = 510.2
$ws.Range("J11").Value = 71
$ws.Range("K11").Value = 510.2
$ws.Range("L11").Value = 71
$ws.Range("M11").Value = -370.2
$ws.Range("N11").Value = -351

$ws.Range("H86").Value = 70111.484
$ws.Range("I86").Value = 1212.7826
$ws.Range("K86").Value = 1212.7826
$ws.Range("M86").Value = -89.7826

$ws.Range("H89").Value = 70111.484
$ws.Range("I89").Value = 1212.7826
$ws.Range("K89").Value = 6063.913
$ws.Range("M89").Value = -447.9130000000005

$ws.Range("H99").Value = 1452.75
$ws.Range("I99").Value = 1458.0454
$ws.Range("J99").Value = 1433.3334
$ws.Range("K99").Value = 1458.0454
$ws.Range("L99").Value = 1433.3334
$ws.Range("M99").Value = 39.95460000000003
$ws.Range("N99").Value = -4429.3334

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H134").Value = 5157.3335
$ws.Range("I134").Value = 3183.52
$ws.Range("K134").Value = 9550.559999999999
$ws.Range("M134").Value = -7015.559999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3801.4285
$ws.Range("J4").Value = 5022
$ws.Range("L4").Value = 5022
$ws.Range("N4").Value = -5246

$ws.Range("H62").Value = 4334.6665
$ws.Range("J62").Value = 4500
$ws.Range("L62").Value = 4500
$ws.Range("N62").Value = -5748

$ws.Range("H65").Value = 4334.6665
$ws.Range("J65").Value = 4500
$ws.Range("L65").Value = 22500
$ws.Range("N65").Value = -28740

$ws.Range("H132").Value = 40281.8
$ws.Range("I132").Value = 3072.0908
$ws.Range("K132").Value = 9216.2724
$ws.Range("M132").Value = -6686.2724

$ws.Range("H137").Value = 109500
$ws.Range("J137").Value = 109500
$ws.Range("L137").Value = 109500
$ws.Range("N137").Value = -119700

$ws.Range("H141").Value = 380326.88
$ws.Range("J141").Value = 406101.44
$ws.Range("L141").Value = 406101.44
$ws.Range("N141").Value = -416461.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 42384396
$ws.Range("I4").Value = 48438596
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 145315788
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -145315676
$ws.Range("N4").Value = -15224

$ws.Range("H23").Value = 197.11111
$ws.Range("J23").Value = 233.16667
$ws.Range("L23").Value = 699.50001
$ws.Range("N23").Value = -1169.50001

$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H117").Value = 2395.7
$ws.Range("J117").Value = 3656.3333
$ws.Range("L117").Value = 10968.9999
$ws.Range("N117").Value = -17852.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 24815.166
$ws.Range("J49").Value = 24815.166
$ws.Range("L49").Value = 24815.166
$ws.Range("N49").Value = -25183.166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9998
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 9998
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 9998
$ws.Range("N2").Value = -10222
$ws.Range("M2").ClearContents()

$ws.Range("H55").Value = 172.75
$ws.Range("I55").Value = 141.9
$ws.Range("J55").Value = 224.16667
$ws.Range("K55").Value = 141.9
$ws.Range("L55").Value = 224.16667
$ws.Range("M55").Value = 31.09999999999999
$ws.Range("N55").Value = -570.1666700000001

$ws.Range("H96").Value = 49000
$ws.Range("J96").Value = 49000
$ws.Range("L96").Value = 49000
$ws.Range("N96").Value = -54492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4200.6
$ws.Range("J2").Value = 4999
$ws.Range("L2").Value = 4999
$ws.Range("N2").Value = -5223

$ws.Range("H81").Value = 12480.883
$ws.Range("J81").Value = 15934.946
$ws.Range("L81").Value = 31869.892
$ws.Range("N81").Value = -33991.892

$ws.Range("H84").Value = 12480.883
$ws.Range("J84").Value = 15934.946
$ws.Range("L84").Value = 159349.46
$ws.Range("N84").Value = -169957.46

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("N130").ClearContents()

$ws.Range("H131").Value = 160730
$ws.Range("J131").Value = 180750
$ws.Range("L131").Value = 180750
$ws.Range("N131").Value = -190830
